$wb = $excel.ActiveWorkbook

$wsCase = $wb.Worksheets.Item("Test Case")

# Copy the formatting of row 2 to row 3 first, so the new row matches styling.
$wsCase.Range("A2:H2").Copy()
$wsCase.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New shared strings must be introduced in this exact order so that the
# resulting sharedStrings.xml table matches the expected layout:
#   ... TC02, "Verify Share to Contact menu at footer- Product Module",
#   "Verify Share to Contact menu at footer- Catalog Module", <Catalog steps>
$wsCase.Range("D3").Value = "TC02"

$wsCase.Range("E2").Value = "Verify Share to Contact menu at footer- Product Module"

$wsCase.Range("E3").Value = "Verify Share to Contact menu at footer- Catalog Module"

$wsCase.Range("H3").Value = "1. Log in to ""http://blubox.shoppinpal.com/s eller/"" website
2. Navigate to Catalog module. Verify the following:
a. User is able to see thumbnails in the selected layout
b. User is able to see checkbox over upper left corner of thumbnails
c. User is able to see options tag over upper roght corner of thumbnails
d. Various details for thumbnails are present"

# Fill in the remaining (already-existing) shared string values for row 3.
$wsCase.Range("A3").Value = 1
$wsCase.Range("B3").Value = $wsCase.Range("B2").Value()
$wsCase.Range("C3").Value = $wsCase.Range("C2").Value()
$wsCase.Range("F3").Value = $wsCase.Range("F2").Value()
$wsCase.Range("G3").Value = $wsCase.Range("G2").Value()

$wsCase.Rows.Item(3).RowHeight = $wsCase.Rows.Item(2).RowHeight

$wsCase.Range("H3").Select()

$ws = $wb.Worksheets.Item("Test Case")
$ws.Activate()
